$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Date" column (F) values for the first two data rows (Supplier, Plant)
# to new refreshed-as-of timestamps.
$ws.Range("F2").Value = "Feb 17, 2022 (02:54:36 EST)"
$ws.Range("F3").Value = "Feb 17, 2022 (02:54:40 EST)"
